$d = $word.ActiveDocument

# Replace every exact, case-sensitive occurrence of $findText with
# $replaceText by locating it with Find (no in-place replace) and then
# assigning Range.Text directly. Doing the substitution this way (and
# not via Find.Execute's own Replace argument) avoids the "smart quote"
# autocorrection that Word applies to typed apostrophes during
# Find&Replace, and keeps unrelated runs elsewhere in the document
# untouched.
function Replace-AllExact($findText, $replaceText, $matchWholeWord) {
    $searchStart = 0
    while ($true) {
        $rng = $d.Range($searchStart, $d.Content.End)
        $found = $rng.Find.Execute($findText, $true, $matchWholeWord, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $rng.Text = $replaceText
        $searchStart = $rng.Start + $replaceText.Length
    }
}

# 1) Salutation/title field result: MR -> MS (both occurrences)
Replace-AllExact "MR" "MS" $true

# 2) Full name field result
Replace-AllExact "INOCENCIO M.  ANGCAYA" "PURISIMA C.  DUNGO" $false

# 3) Position field result
Replace-AllExact "Casual Employee" "Ticket Checker" $false

# 4) Office field result
Replace-AllExact "Vice Mayor's Office Detailed At Civil Security Unit" "City Treasurer's Office" $false

# 5) Last day of service field result
Replace-AllExact "July 29, 2021" "December 31, 2022" $false

# 6) Vacation leave amount
Replace-AllExact "  38.708" "  22.857" $false

# 7) Sick leave amount
Replace-AllExact "  53.708" "  57.250" $false

# 8) Total leave amount
Replace-AllExact "  92.416" "  80.107" $false

# 9) Requester gender title: Mr -> Ms (case sensitive whole word so "MR" stays intact)
Replace-AllExact "Mr" "Ms" $true

# 10) Requester last name
Replace-AllExact "Angcaya" "Dungo" $false

# 11) "Issued this 22nd day of December" -> "Issued this 23rd day of December"
# Replicates Word splitting "22" into two runs ("2" then a freshly
# typed "3") while leaving the following superscript "nd" -> "rd".
$positions = @()
$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("Issued this 22", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $positions += $rng.End
    $searchStart = $rng.End
}

foreach ($endPos in $positions) {
    $secondDigit = $d.Range($endPos - 1, $endPos)
    $secondDigit.Text = ""
    $secondDigit.InsertAfter("3")
    # Force the engine to keep the newly typed "3" as its own run
    # instead of silently re-merging it with the preceding "2".
    $secondDigit.Font.Bold = 1
    $secondDigit.Font.Bold = 0
}

Replace-AllExact "nd" "rd" $false
